$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
# Row 74
$ws.Range("H74").Value = 103356.44
$ws.Range("I74").Value = 152770.4
$ws.Range("J74").Value = 20999.834
$ws.Range("K74").Value = 152770.4
$ws.Range("L74").Value = 20999.834
$ws.Range("M74").Value = -151834.4
$ws.Range("N74").Value = -22871.834
# Row 77
$ws.Range("H77").Value = 103356.44
$ws.Range("I77").Value = 152770.4
$ws.Range("J77").Value = 20999.834
$ws.Range("K77").Value = 763852
$ws.Range("L77").Value = 104999.17
$ws.Range("M77").Value = -759172
$ws.Range("N77").Value = -114359.17
# Row 80
$ws.Range("H80").Value = 386.7857
$ws.Range("I80").Value = 294
$ws.Range("J80").Value = 479.57144
$ws.Range("K80").Value = 882
$ws.Range("L80").Value = 1438.71432
$ws.Range("M80").Value = 116
$ws.Range("N80").Value = -3434.71432
# Row 83
$ws.Range("H83").Value = 386.7857
$ws.Range("I83").Value = 294
$ws.Range("J83").Value = 479.57144
$ws.Range("K83").Value = 2646
$ws.Range("L83").Value = 4316.14296
$ws.Range("M83").Value = 2346
$ws.Range("N83").Value = -14300.14296
# Row 116
$ws.Range("H116").Value = 8977.777
$ws.Range("I116").Value = 10600.6
$ws.Range("J116").Value = 6949.25
$ws.Range("K116").Value = 10600.6
$ws.Range("L116").Value = 6949.25
$ws.Range("M116").Value = -7158.6
$ws.Range("N116").Value = -13833.25

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = 0
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("N39").Value = 0
# Row 102
$ws.Range("H102").Value = 3884.6667
$ws.Range("I102").Value = 3661.6
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 3661.6
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -2039.6
$ws.Range("N102").Value = -8244
# Row 132
$ws.Range("H132").Value = 1975
$ws.Range("I132").Value = 1975
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5925
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3395

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
# Row 19
$ws.Range("H19").Value = 33342.668
$ws.Range("I19").Value = 9
$ws.Range("J19").Value = 100010
$ws.Range("K19").Value = 9
$ws.Range("L19").Value = 100010
$ws.Range("M19").Value = 164
$ws.Range("N19").Value = -100356
# Row 99
$ws.Range("H99").Value = 1214.9
$ws.Range("I99").Value = 1322.1111
$ws.Range("J99").Value = 250
$ws.Range("K99").Value = 1322.1111
$ws.Range("L99").Value = 250
$ws.Range("M99").Value = 175.8888999999999
$ws.Range("N99").Value = -3246
# Row 105
$ws.Range("H105").Value = 2233.1667
$ws.Range("I105").Value = 2179.8
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2179.8
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -432.8000000000002
$ws.Range("N105").Value = -5994
# Row 107
$ws.Range("H107").Value = 1737.2142
$ws.Range("I107").Value = 1756.7273
$ws.Range("J107").Value = 1665.6666
$ws.Range("K107").Value = 1756.7273
$ws.Range("L107").Value = 1665.6666
$ws.Range("M107").Value = 163.2727
$ws.Range("N107").Value = -5505.6666

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
# Row 58
$ws.Range("H58").Value = 1851.75
$ws.Range("I58").Value = 1640.1428
$ws.Range("J58").Value = 3333
$ws.Range("K58").Value = 1640.1428
$ws.Range("L58").Value = 3333
$ws.Range("M58").Value = -1437.1428
$ws.Range("N58").Value = -3739
# Row 105
$ws.Range("H105").Value = 3357.1428
$ws.Range("I105").Value = 2500
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 2500
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -753
$ws.Range("N105").Value = -7994
# Row 136
$ws.Range("H136").Value = 1851.75
$ws.Range("I136").Value = 1640.1428
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 4920.428400000001
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -2370.428400000001
$ws.Range("N136").Value = -15099

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
# Row 113
$ws.Range("H113").Value = 1637
$ws.Range("I113").Value = 1156.8
$ws.Range("J113").Value = 2237.25
$ws.Range("K113").Value = 3470.4
$ws.Range("L113").Value = 6711.75
$ws.Range("M113").Value = -1300.4
$ws.Range("N113").Value = -11051.75
# Row 132
$ws.Range("H132").Value = 4699.4
$ws.Range("I132").Value = 4374.25
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 39368.25
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -36838.25
$ws.Range("N132").Value = -59060

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
# Row 22
$ws.Range("H22").Value = 5666.6665
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 5666.6665
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 5666.6665
$ws.Range("N22").Value = -6724.6665
# Row 70
$ws.Range("H70").Value = 3253.75
$ws.Range("I70").Value = 2672
$ws.Range("J70").Value = 4999
$ws.Range("K70").Value = 2672
$ws.Range("L70").Value = 4999
$ws.Range("M70").Value = -2402
$ws.Range("N70").Value = -5539
# Row 73
$ws.Range("H73").Value = 3253.75
$ws.Range("I73").Value = 2672
$ws.Range("J73").Value = 4999
$ws.Range("K73").Value = 2672
$ws.Range("L73").Value = 4999
$ws.Range("M73").Value = -1736
$ws.Range("N73").Value = -6871
# Row 132
$ws.Range("H132").Value = 15933.333
$ws.Range("I132").Value = 15933.333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 47799.999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -45269.999

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
# Row 13
$ws.Range("H13").Value = 411.1111
$ws.Range("I13").Value = 150
$ws.Range("J13").Value = 1325
$ws.Range("K13").Value = 150
$ws.Range("L13").Value = 1325
$ws.Range("M13").Value = -10
$ws.Range("N13").Value = -1605
# Row 63
$ws.Range("H63").Value = 87077
$ws.Range("I63").Value = 87077
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 87077
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -86328
# Row 66
$ws.Range("H66").Value = 87077
$ws.Range("I66").Value = 87077
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 261231
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -257487
# Row 93
$ws.Range("H93").Value = 1274.6875
$ws.Range("I93").Value = 1105.6
$ws.Range("J93").Value = 1556.5
$ws.Range("K93").Value = 1105.6
$ws.Range("L93").Value = 1556.5
$ws.Range("M93").Value = 142.4000000000001
$ws.Range("N93").Value = -4052.5
# Row 100
$ws.Range("H100").Value = 2818.7222
$ws.Range("I100").Value = 2722.6667
$ws.Range("J100").Value = 3299
$ws.Range("K100").Value = 2722.6667
$ws.Range("L100").Value = 3299
$ws.Range("M100").Value = -2181.6667
$ws.Range("N100").Value = -4381
# Row 136
$ws.Range("H136").Value = 1519.9
$ws.Range("I136").Value = 1100.8572
$ws.Range("J136").Value = 2497.6667
$ws.Range("K136").Value = 3302.5716
$ws.Range("L136").Value = 7493.000100000001
$ws.Range("M136").Value = -752.5715999999998
$ws.Range("N136").Value = -12593.0001

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
# Row 74
$ws.Range("H74").Value = 22337.666
$ws.Range("I74").Value = 20627.5
$ws.Range("J74").Value = 23192.75
$ws.Range("K74").Value = 20627.5
$ws.Range("L74").Value = 23192.75
$ws.Range("M74").Value = -19691.5
$ws.Range("N74").Value = -25064.75
# Row 77
$ws.Range("H77").Value = 22337.666
$ws.Range("I77").Value = 20627.5
$ws.Range("J77").Value = 23192.75
$ws.Range("K77").Value = 61882.5
$ws.Range("L77").Value = 69578.25
$ws.Range("M77").Value = -57202.5
$ws.Range("N77").Value = -78938.25
# Row 81
$ws.Range("H81").Value = 2507250.2
$ws.Range("I81").Value = 9666.333000000001
$ws.Range("J81").Value = 10000002
$ws.Range("K81").Value = 19332.666
$ws.Range("L81").Value = 20000004
$ws.Range("M81").Value = -18271.666
$ws.Range("N81").Value = -20002126
# Row 84
$ws.Range("H84").Value = 2507250.2
$ws.Range("I84").Value = 9666.333000000001
$ws.Range("J84").Value = 10000002
$ws.Range("K84").Value = 96663.33
$ws.Range("L84").Value = 100000020
$ws.Range("M84").Value = -91359.33
$ws.Range("N84").Value = -100010628
